# ---------------------------------------------------------------------------
# Apply the "TidyBuddy" sample-data refresh:
#   - primer_index: the two 96-well "primer name" blocks (rows 2-9, rows 12-19)
#     get re-labelled from the old `5_A01_set323`.. style names to generic
#     `ind_0001`.. placeholders, and the second block's content is replaced
#     with the same repeating ind_0001..ind_0012 sequence used elsewhere.
#     The ad-hoc per-cell formatting that used to sit on those two blocks is
#     cleared, and ends up instead on a new, otherwise-empty, row 21.
#   - Replicate: the second 8-row block (rows 12-19) changes from a repeating
#     1..6 pattern to a constant-per-row value (1,2,3,...,8).
#   - Selections on all four sheets move to reflect where the user left the
#     cursor afterwards.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- primer_index -----------------------------------------------------------
$ws1 = $wb.Worksheets.Item("primer_index")
$ws1.Activate()

# Row 2-9, columns B..M: same layout as before (6 cells of "ind_000X" then
# 6 cells of "ind_000(X+1)"), just renamed from the old 5_xxx_setNNN strings.
$block1 = @(
    @("ind_0001","ind_0001","ind_0001","ind_0001","ind_0001","ind_0001","ind_0002","ind_0002","ind_0002","ind_0002","ind_0002","ind_0002"),
    @("ind_0003","ind_0003","ind_0003","ind_0003","ind_0003","ind_0003","ind_0004","ind_0004","ind_0004","ind_0004","ind_0004","ind_0004"),
    @("ind_0005","ind_0005","ind_0005","ind_0005","ind_0005","ind_0005","ind_0006","ind_0006","ind_0006","ind_0006","ind_0006","ind_0006"),
    @("ind_0007","ind_0007","ind_0007","ind_0007","ind_0007","ind_0007","ind_0008","ind_0008","ind_0008","ind_0008","ind_0008","ind_0008"),
    @("ind_0009","ind_0009","ind_0009","ind_0009","ind_0009","ind_0009","ind_0010","ind_0010","ind_0010","ind_0010","ind_0010","ind_0010"),
    @("ind_0011","ind_0011","ind_0011","ind_0011","ind_0011","ind_0011","ind_0012","ind_0012","ind_0012","ind_0012","ind_0012","ind_0012"),
    @("ind_0013","ind_0013","ind_0013","ind_0013","ind_0013","ind_0013","ind_0014","ind_0014","ind_0014","ind_0014","ind_0014","ind_0014"),
    @("ind_0015","ind_0015","ind_0015","ind_0015","ind_0015","ind_0015","ind_0016","ind_0016","ind_0016","ind_0016","ind_0016","ind_0016")
)

for ($i = 0; $i -lt $block1.Length; $i++) {
    $r = 2 + $i
    $rowVals = $block1[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $c = 2 + $j
        $ws1.Cells.Item($r, $c).Value = $rowVals[$j]
    }
}

# Row 12-19, columns B..M: every row now repeats the same ind_0001..ind_0012
# sequence (instead of each row having its own pair of set-names).
$rowPattern = @("ind_0001","ind_0002","ind_0003","ind_0004","ind_0005","ind_0006","ind_0007","ind_0008","ind_0009","ind_0010","ind_0011","ind_0012")

for ($r = 12; $r -le 19; $r++) {
    for ($j = 0; $j -lt $rowPattern.Length; $j++) {
        $c = 2 + $j
        $ws1.Cells.Item($r, $c).Value = $rowPattern[$j]
    }
}

# The old per-cell style ("s=1", a Calibri 11pt font override) is removed
# from those two data blocks...
$ws1.Range("B2:M9").ClearFormats()
$ws1.Range("B12:M19").ClearFormats()

# ...and shows up instead on a new, empty row 21 (columns B..Q), copied from
# an existing cell that still carries that same style.
$ws1.Range("N2").Copy()
$ws1.Range("B21:Q21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws1.Range("M26").Select()

# --- Replicate ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Replicate")
$ws2.Activate()

for ($i = 0; $i -lt 8; $i++) {
    $r = 12 + $i
    $val = $i + 1
    for ($c = 2; $c -le 13; $c++) {
        $ws2.Cells.Item($r, $c).Value = $val
    }
}

$ws2.Range("D24").Select()

# --- PrimerPlate --------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("PrimerPlate")
$ws3.Activate()
$ws3.Range("A10:XFD23").Select()

# --- PrimerWell ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("PrimerWell")
$ws4.Activate()
$ws4.Range("G23").Select()

# Leave the original tab (primer_index) active again, matching the saved file.
$ws1.Activate()
